# Auto-applied scheduled-runner update to price/profit figures across
# the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR Leve-profit tables.
# For each touched row, refresh the market-price-derived columns
# (currentAveragePrice[NQ/HQ], LevePrice[NQ/HQ], LeveProfit[NQ/HQ]).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 769163
$ws.Range("I17").Value = 243.75
$ws.Range("K17").Value = 731.25
$ws.Range("M17").Value = -563.25

$ws.Range("H62").Value = 3980422.2
$ws.Range("I62").Value = 6181601
$ws.Range("J62").Value = 18300
$ws.Range("K62").Value = 6181601
$ws.Range("L62").Value = 18300
$ws.Range("M62").Value = -6180977
$ws.Range("N62").Value = -19548

$ws.Range("H65").Value = 3980422.2
$ws.Range("I65").Value = 6181601
$ws.Range("J65").Value = 18300
$ws.Range("K65").Value = 30908005
$ws.Range("L65").Value = 91500
$ws.Range("M65").Value = -30904885
$ws.Range("N65").Value = -97740

$ws.Range("H86").Value = 1760
$ws.Range("I86").Value = 1000
$ws.Range("J86").Value = 1950
$ws.Range("K86").Value = 1000
$ws.Range("L86").Value = 1950
$ws.Range("M86").Value = 123
$ws.Range("N86").Value = -4196

$ws.Range("H89").Value = 1760
$ws.Range("I89").Value = 1000
$ws.Range("J89").Value = 1950
$ws.Range("K89").Value = 5000
$ws.Range("L89").Value = 9750
$ws.Range("M89").Value = 616
$ws.Range("N89").Value = -20982

$ws.Range("H107").Value = 427762.88
$ws.Range("I107").Value = 555971.9399999999
$ws.Range("J107").Value = 399.33334
$ws.Range("K107").Value = 555971.9399999999
$ws.Range("L107").Value = 399.33334
$ws.Range("M107").Value = -554051.9399999999
$ws.Range("N107").Value = -4239.33334

$ws.Range("H113").Value = 146857.86
$ws.Range("I113").Value = 170500.83
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 170500.83
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = -167246.83
$ws.Range("N113").Value = -11508

$ws.Range("H132").Value = 190826.73
$ws.Range("I132").Value = 203365.53
$ws.Range("J132").Value = 40361.2
$ws.Range("K132").Value = 610096.59
$ws.Range("L132").Value = 121083.6
$ws.Range("M132").Value = -607566.59
$ws.Range("N132").Value = -126143.6

$ws.Range("H135").Value = 1317.1077
$ws.Range("I135").Value = 1169.5834
$ws.Range("J135").Value = 1733.6471
$ws.Range("K135").Value = 10526.2506
$ws.Range("L135").Value = 15602.8239
$ws.Range("M135").Value = -7991.250599999999
$ws.Range("N135").Value = -20672.8239

$ws.Range("H137").Value = 21740220
$ws.Range("I137").Value = 29412550
$ws.Range("J137").Value = 1955.0834
$ws.Range("K137").Value = 88237650
$ws.Range("L137").Value = 5865.2502
$ws.Range("M137").Value = -88235100
$ws.Range("N137").Value = -10965.2502

$ws.Range("H138").Value = 7852384
$ws.Range("I138").Value = 946997.0600000001
$ws.Range("J138").Value = 166676290
$ws.Range("K138").Value = 2840991.18
$ws.Range("L138").Value = 500028870
$ws.Range("M138").Value = -2835851.18
$ws.Range("N138").Value = -500039150

$ws.Range("H141").Value = 2179.25
$ws.Range("I141").Value = 1270.2693
$ws.Range("J141").Value = 6118.1665
$ws.Range("K141").Value = 3810.8079
$ws.Range("L141").Value = 18354.4995
$ws.Range("M141").Value = 1369.1921
$ws.Range("N141").Value = -28714.4995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 125424.875
$ws.Range("I5").Value = 167099.83
$ws.Range("J5").Value = 400
$ws.Range("K5").Value = 167099.83
$ws.Range("L5").Value = 400
$ws.Range("M5").Value = -166987.83
$ws.Range("N5").Value = -624

$ws.Range("H74").Value = 3426.6667
$ws.Range("I74").Value = 1068
$ws.Range("J74").Value = 11409.846
$ws.Range("K74").Value = 1068
$ws.Range("L74").Value = 11409.846
$ws.Range("M74").Value = -194
$ws.Range("N74").Value = -13157.846

$ws.Range("H77").Value = 3426.6667
$ws.Range("I77").Value = 1068
$ws.Range("J77").Value = 11409.846
$ws.Range("K77").Value = 5340
$ws.Range("L77").Value = 57049.23
$ws.Range("M77").Value = -972
$ws.Range("N77").Value = -65785.23

$ws.Range("H110").Value = 790.75
$ws.Range("I110").Value = 611.5833
$ws.Range("J110").Value = 1328.25
$ws.Range("K110").Value = 611.5833
$ws.Range("L110").Value = 1328.25
$ws.Range("M110").Value = 1433.4167
$ws.Range("N110").Value = -5418.25

$ws.Range("H132").Value = 1744.3906
$ws.Range("I132").Value = 1527.7407
$ws.Range("K132").Value = 4583.2221
$ws.Range("M132").Value = -2053.2221

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 125424.875
$ws.Range("I4").Value = 167099.83
$ws.Range("J4").Value = 400
$ws.Range("K4").Value = 167099.83
$ws.Range("L4").Value = 400
$ws.Range("M4").Value = -166984.83
$ws.Range("N4").Value = -630

$ws.Range("H20").Value = 3999.75
$ws.Range("I20").Value = 3999.6667
$ws.Range("K20").Value = 3999.6667
$ws.Range("M20").Value = -3752.6667

$ws.Range("H94").Value = 691.2857
$ws.Range("I94").Value = 522.03125
$ws.Range("K94").Value = 522.03125
$ws.Range("M94").Value = -71.03125

$ws.Range("H107").Value = 889.4
$ws.Range("I107").Value = 902.93335
$ws.Range("J107").Value = 848.8
$ws.Range("K107").Value = 902.93335
$ws.Range("L107").Value = 848.8
$ws.Range("M107").Value = 1017.06665
$ws.Range("N107").Value = -4688.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 917.8
$ws.Range("I16").Value = 428.2857
$ws.Range("K16").Value = 428.2857
$ws.Range("M16").Value = -141.2857

$ws.Range("H31").Value = 1835.0714
$ws.Range("I31").Value = 1033.909
$ws.Range("K31").Value = 1033.909
$ws.Range("M31").Value = -738.9090000000001

$ws.Range("H34").Value = 1835.0714
$ws.Range("I34").Value = 1033.909
$ws.Range("K34").Value = 1033.909
$ws.Range("M34").Value = -831.9090000000001

$ws.Range("H58").Value = 2344.75
$ws.Range("I58").Value = 982.7692
$ws.Range("J58").Value = 4874.143
$ws.Range("K58").Value = 982.7692
$ws.Range("L58").Value = 4874.143
$ws.Range("M58").Value = -779.7692
$ws.Range("N58").Value = -5280.143

$ws.Range("H107").Value = 304.57144
$ws.Range("I107").Value = 186.46666
$ws.Range("J107").Value = 599.8333
$ws.Range("K107").Value = 186.46666
$ws.Range("L107").Value = 599.8333
$ws.Range("M107").Value = 1733.53334
$ws.Range("N107").Value = -4439.8333

$ws.Range("H113").Value = 917.8
$ws.Range("I113").Value = 428.2857
$ws.Range("K113").Value = 428.2857
$ws.Range("M113").Value = 1741.7143

$ws.Range("H135").Value = 41236.668
$ws.Range("J135").Value = 41236.668
$ws.Range("L135").Value = 41236.668
$ws.Range("N135").Value = -51376.668

$ws.Range("H136").Value = 2344.75
$ws.Range("I136").Value = 982.7692
$ws.Range("J136").Value = 4874.143
$ws.Range("K136").Value = 2948.3076
$ws.Range("L136").Value = 14622.429
$ws.Range("M136").Value = -398.3076000000001
$ws.Range("N136").Value = -19722.429

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 29.291666
$ws.Range("I12").Value = 19.307692
$ws.Range("J12").Value = 41.090908
$ws.Range("K12").Value = 57.92307599999999
$ws.Range("L12").Value = 123.272724
$ws.Range("M12").Value = 115.076924
$ws.Range("N12").Value = -469.272724

$ws.Range("H131").Value = 2216.3333
$ws.Range("J131").Value = 2331.5806
$ws.Range("L131").Value = 6994.7418
$ws.Range("N131").Value = -17074.7418

$ws.Range("H132").Value = 1916.1428
$ws.Range("J132").Value = 1978.9231
$ws.Range("L132").Value = 17810.3079
$ws.Range("N132").Value = -22870.3079

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 54.833332
$ws.Range("I2").Value = 16.333334
$ws.Range("K2").Value = 16.333334
$ws.Range("M2").Value = 96.66666599999999

$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("N103").ClearContents()

$ws.Range("H126").Value = 2180.4614
$ws.Range("I126").Value = 1659.2
$ws.Range("J126").Value = 2506.25
$ws.Range("K126").Value = 4977.6
$ws.Range("L126").Value = 7518.75
$ws.Range("M126").Value = -2507.6
$ws.Range("N126").Value = -12458.75

$ws.Range("H132").Value = 2549.7896
$ws.Range("I132").Value = 2217.2292
$ws.Range("J132").Value = 4323.4443
$ws.Range("K132").Value = 6651.687600000001
$ws.Range("L132").Value = 12970.3329
$ws.Range("M132").Value = -4121.687600000001
$ws.Range("N132").Value = -18030.3329

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 5001491.5
$ws.Range("I16").Value = 5264675.5
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 5264675.5
$ws.Range("L16").Value = 1000
$ws.Range("M16").Value = -5264505.5
$ws.Range("N16").Value = -1340

$ws.Range("H22").Value = 13011.111
$ws.Range("I22").Value = 1700
$ws.Range("J22").Value = 16242.857
$ws.Range("K22").Value = 1700
$ws.Range("L22").Value = 16242.857
$ws.Range("M22").Value = -1405
$ws.Range("N22").Value = -16832.857

$ws.Range("H27").Value = 13011.111
$ws.Range("I27").Value = 1700
$ws.Range("J27").Value = 16242.857
$ws.Range("K27").Value = 1700
$ws.Range("L27").Value = 16242.857
$ws.Range("M27").Value = -1593
$ws.Range("N27").Value = -16456.857

$ws.Range("H43").Value = 10606.533
$ws.Range("J43").Value = 7792.7144
$ws.Range("L43").Value = 7792.7144
$ws.Range("N43").Value = -8178.7144

$ws.Range("H46").Value = 2822.8572
$ws.Range("I46").Value = 2500
$ws.Range("K46").Value = 2500
$ws.Range("M46").Value = -2312

$ws.Range("H68").Value = 2349.8333
$ws.Range("I68").Value = 1750
$ws.Range("J68").Value = 2649.75
$ws.Range("K68").Value = 1750
$ws.Range("L68").Value = 2649.75
$ws.Range("M68").Value = -1001
$ws.Range("N68").Value = -4147.75

$ws.Range("H71").Value = 2349.8333
$ws.Range("I71").Value = 1750
$ws.Range("J71").Value = 2649.75
$ws.Range("K71").Value = 8750
$ws.Range("L71").Value = 13248.75
$ws.Range("M71").Value = -5006
$ws.Range("N71").Value = -20736.75

$ws.Range("H93").Value = 1600.8
$ws.Range("I93").Value = 1366.6666
$ws.Range("J93").Value = 1952
$ws.Range("K93").Value = 1366.6666
$ws.Range("L93").Value = 1952
$ws.Range("M93").Value = -118.6666
$ws.Range("N93").Value = -4448

$ws.Range("H132").Value = 4085.4795
$ws.Range("I132").Value = 4106.3335
$ws.Range("J132").Value = 4026.2104
$ws.Range("K132").Value = 12319.0005
$ws.Range("L132").Value = 12078.6312
$ws.Range("M132").Value = -9789.000499999998
$ws.Range("N132").Value = -17138.6312

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 15492.928
$ws.Range("I136").Value = 16757.111
$ws.Range("J136").Value = 2219
$ws.Range("K136").Value = 50271.333
$ws.Range("L136").Value = 6657
$ws.Range("M136").Value = -47721.333
$ws.Range("N136").Value = -11757
